$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values replacing old ones, per row (2..15)
$newValues = @{
    2  = 8
    3  = 11
    4  = 5
    5  = 13
    6  = 10
    7  = 8
    8  = 9
    9  = 9
    10 = 9
    11 = 10
    12 = 12
    13 = 6
    14 = 5
    15 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
